$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Alternative")
$ws.Range("I2").Value = 9.77
$ws.Range("J2").Value = 9.77
$ws.Range("K2").Value = 12.59
$ws.Range("M2").Value = -11.15
$ws.Range("I3").Value = 1.36
$ws.Range("J3").Value = 9.19
$ws.Range("K3").Value = 8.24
$ws.Range("L3").Value = 7.11
$ws.Range("M3").Value = -4.08
$ws.Range("I4").Value = 1.23
$ws.Range("J4").Value = 3.14
$ws.Range("K4").Value = 7.41
$ws.Range("M4").Value = -8.880000000000001
$ws.Range("I5").Value = -8.41
$ws.Range("J5").Value = -8.41
$ws.Range("K5").Value = 45.78
$ws.Range("L5").Value = 26.82
$ws.Range("M5").Value = -83.40000000000001
$ws.Range("I6").Value = -0.16
$ws.Range("J6").Value = 2.22
$ws.Range("K6").Value = 6.72
$ws.Range("L6").Value = 7.94
$ws.Range("M6").Value = -9.17
$ws.Range("N6").Value = -0.45
$ws.Range("K7").Value = 5.66
$ws.Range("M7").Value = -7.72
$ws.Range("N7").Value = -0.64

$ws = $wb.Worksheets.Item("Bond")
$ws.Range("I2").Value = -3.76
$ws.Range("M2").Value = -6.47
$ws.Range("N3").Value = -2.84
$ws.Range("I4").Value = -2.23
$ws.Range("J4").Value = 3.2
$ws.Range("L4").Value = 15.65
$ws.Range("M4").Value = -26.48
$ws.Range("M5").Value = -16.53

$ws = $wb.Worksheets.Item("Equity")
$ws.Range("I2").Value = 12.88
$ws.Range("J2").Value = 14.53
$ws.Range("K2").Value = 20.02
$ws.Range("M2").Value = -18.46
$ws.Range("I3").Value = 8.66
$ws.Range("J3").Value = 11.87
$ws.Range("M3").Value = -15.6
$ws.Range("I4").Value = 4.18
$ws.Range("J4").Value = 6.44
$ws.Range("K4").Value = 14.68
$ws.Range("M4").Value = -17.9
$ws.Range("M5").Value = -9.83
$ws.Range("I6").Value = 1.37
$ws.Range("J6").Value = 5.23
$ws.Range("K6").Value = 26.91
$ws.Range("M6").Value = -39.16
$ws.Range("N6").Value = 0
$ws.Range("I7").Value = 5.37
$ws.Range("K7").Value = 15.64
$ws.Range("L7").Value = 10.7
$ws.Range("M7").Value = -16.91
$ws.Range("N7").Value = 0.23
$ws.Range("I8").Value = 1.52
$ws.Range("J8").Value = 4.94
$ws.Range("M8").Value = -18.98
$ws.Range("I9").Value = -6.32
$ws.Range("J9").Value = -0.86
$ws.Range("L9").Value = 16.21
$ws.Range("M9").Value = -23.3
$ws.Range("I10").Value = 2.46
$ws.Range("J10").Value = 7.52
$ws.Range("K10").Value = 26.3
$ws.Range("L10").Value = 16.48
$ws.Range("M10").Value = -35.98
$ws.Range("I11").Value = -1.13
$ws.Range("J11").Value = 11.8
$ws.Range("L11").Value = 9.050000000000001
$ws.Range("M11").Value = -0.97

Write-Host "Applied all updates"
